# Updated symbol list on Mon Jan  9 21:55:04 UTC 2023 with GitHub Actions
#
# Refreshes the Price (col D) and Volume(1h) (col E) snapshot values for
# the crypto rows on the active sheet. Source cells are stored as plain
# text (e.g. "273.54", "1.26%"), so each write forces a text number
# format before assigning the new literal string and then resets the
# cell style back to Normal so no stray formatting/number-format is
# introduced versus the original file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "272.87"
Set-TextCell "E2" "1.04%"
Set-TextCell "D3" "26.83"
Set-TextCell "E3" "0.40%"
Set-TextCell "D4" "4.908"
Set-TextCell "E4" "3.94%"
Set-TextCell "D5" "0.06317"
Set-TextCell "E5" "3.23%"
Set-TextCell "D6" "6.918"
Set-TextCell "E6" "2.56%"
Set-TextCell "D7" "3.349"
Set-TextCell "E7" "5.44%"
Set-TextCell "D8" "1.330"
Set-TextCell "E8" "48.72%"
Set-TextCell "D9" "0.8851"
Set-TextCell "E9" "3.28%"
Set-TextCell "E10" "2.79%"
Set-TextCell "D11" "0.05095"
Set-TextCell "E11" "2.09%"
Set-TextCell "D12" "0.07396"
Set-TextCell "E12" "3.94%"
Set-TextCell "D13" "0.03197"
Set-TextCell "E13" "0.41%"
Set-TextCell "D14" "0.09051"
Set-TextCell "E14" "0.17%"
Set-TextCell "D15" "0.001582"
Set-TextCell "E15" "2.97%"
Set-TextCell "D16" "0.0006333"
Set-TextCell "E16" "4.14%"
Set-TextCell "D17" "0.006081"
Set-TextCell "E17" "0.50%"
Set-TextCell "D18" "3.472"
Set-TextCell "E18" "0.26%"
Set-TextCell "D19" "2.284"
Set-TextCell "E19" "0.88%"
Set-TextCell "D20" "0.3142"
Set-TextCell "E20" "1.70%"
Set-TextCell "D21" "0.1333"
Set-TextCell "E21" "4.04%"
Set-TextCell "D22" "3.910"
Set-TextCell "E22" "1.68%"
Set-TextCell "D23" "0.04353"
Set-TextCell "E23" "2.41%"
Set-TextCell "D24" "0.001181"
Set-TextCell "E24" "0.38%"
Set-TextCell "D25" "0.003639"
Set-TextCell "E25" "-12.22%"
Set-TextCell "D26" "0.0001203"
Set-TextCell "E26" "0.30%"
Set-TextCell "E27" "1.18%"
Set-TextCell "D40" "0.04053"
Set-TextCell "E40" "2.59%"
Set-TextCell "D41" "0.006628"
Set-TextCell "E41" "57.88%"
Set-TextCell "D42" "0.1164"
Set-TextCell "E42" "3.95%"
Set-TextCell "D43" "0.002207"
Set-TextCell "E43" "8.31%"
Set-TextCell "D44" "0.01261"
Set-TextCell "E44" "-5.08%"
Set-TextCell "D45" "0.00005353"
Set-TextCell "E45" "3.97%"
Set-TextCell "E46" "149.29%"
Set-TextCell "D47" "0.02126"
Set-TextCell "E47" "-13.14%"
